$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.269.81'
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = '1.444.31'
$ws.Range("E3").Value = '  +1.88%  '
$ws.Range("E4").Value = '  +0.84%  '
$ws.Range("D5").Value = "'0.9267"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -7.23%  '
$ws.Range("D6").Value = "'274.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("D7").Value = "'0.3644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.08%  '
$ws.Range("D8").Value = "'0.3083"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.50%  '
$ws.Range("D9").Value = "'39.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.67%  '
$ws.Range("D10").Value = "'1.025"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.77%  '
$ws.Range("D11").Value = "'0.06531"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D12").Value = "'0.9978"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").Value = "'5.359"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.44%  '
$ws.Range("D14").Value = "'17.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.03%  '
$ws.Range("D15").Value = "'6.071"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.09%  '
$ws.Range("D16").Value = "'0.00001013"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.61%  '
$ws.Range("D17").Value = '1.438.29'
$ws.Range("E17").Value = '  +1.46%  '
$ws.Range("D18").Value = "'0.9439"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.39%  '
$ws.Range("D19").Value = "'0.05663"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.83%  '
$ws.Range("D20").Value = "'69.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.03%  '
$ws.Range("D21").Value = "'5.382"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.18%  '
$ws.Range("D22").Value = "'14.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.46%  '
$ws.Range("E23").Value = '  -3.17%  '
$ws.Range("D24").Value = "'2.250"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.59%  '
$ws.Range("D25").Value = '20.264.56'
$ws.Range("E25").Value = '  +0.97%  '
$ws.Range("D26").Value = "'140.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.50%  '
$ws.Range("D27").Value = "'2.039"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -10.20%  '
$ws.Range("D28").Value = "'17.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.89%  '
$ws.Range("D29").Value = '1.591.25'
$ws.Range("E29").Value = '  +0.77%  '
$ws.Range("D30").Value = "'111.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.06%  '
$ws.Range("D31").Value = "'3.994"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("D32").Value = "'4.826"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -10.13%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = "'0.7838"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.67%  '
$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").Value = "'0.07695"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.17%  '
$ws.Range("D35").Value = "'1.453"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.77%  '
$ws.Range("D36").Value = "'0.05667"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.95%  '
$ws.Range("D37").Value = "'4.656"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.77%  '
$ws.Range("E38").Value = '  +1.72%  '
$ws.Range("D39").Value = "'0.01993"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.86%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = "'10.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.42%  '
$ws.Range("B41").Value = 'Frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D41").Value = "'0.9371"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.04%  '
$ws.Range("D42").Value = "'0.1843"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.77%  '
$ws.Range("D43").Value = "'6.992"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -17.69%  '
$ws.Range("D44").Value = "'0.5189"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.31%  '
$ws.Range("D45").Value = "'3.468"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.79%  '
$ws.Range("D46").Value = "'11.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.30%  '
$ws.Range("D47").Value = "'115.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.88%  '
$ws.Range("D48").Value = "'0.5099"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.28%  '
$ws.Range("D49").Value = "'1.731"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.23%  '
$ws.Range("D50").Value = "'0.06370"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.45%  '
$ws.Range("D51").Value = "'0.9894"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.75%  '
